$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the "Meta description: ..." paragraph that currently sits
#    right after the title (Heading1) paragraph.
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# ------------------------------------------------------------------
# 2. At the end of the document, insert a new bold paragraph carrying
#    the page title text, right before the final (italic) paragraph.
# ------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
$insertionPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Age of the Gods Goddess of Wisdom for Free | Review</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertionPoint.InsertXML($xml)

# InsertXML leaves a trailing blank paragraph used purely to force the
# paragraph break; drop it now that the break has been created.
$blankPara = $d.Paragraphs.Item($lastIndex + 1)
$blankPara.Range.Delete()

# ------------------------------------------------------------------
# 3. Replace the text of the final (italic) paragraph - formerly the
#    image-prompt text - with the meta description text, keeping the
#    paragraph's existing italic formatting intact.
# ------------------------------------------------------------------
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$finalPara.Range.Find.Execute(
    "Create a feature image for `"Age of the Gods: Goddess of Wisdom`". The image should be in a cartoon style and feature a happy Maya warrior with glasses. The warrior should be standing in front of an ancient Greek temple while holding a shield with the game title written on it. Athena, the central figure of the game, should be standing next to the warrior with a confident stance. The symbols of the game, including the Gorgoneion, Olive Branches, Helmets, and Armor, should be seen floating around the two figures. The image should be colorful and eye-catching to attract players' attention.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Find out everything you need to know about Age of the Gods Goddess of Wisdom, a Playtech online slot game with a Greek mythology theme. Play it for free and hit any of the 4 progressive jackpots at any time.",
    2)
